# updata DOCs for wk2 sprint2
#
# Wk1 Sprint2 burndown: 1 story point was completed on day 6 (2019-03-30,
# row 7). Record it in the "Completed" column; the BurnDn/SUM formulas
# already on the sheet pick the change up on recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D7").Value = 1

# Leave the cursor sitting on the cell that was just edited.
$ws.Range("D7").Select()
